# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly generated output, per commit:
#   "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# Row -> New Value mapping for column F, shared between both sheets
$updates = @{
    3  = 528
    4  = 1526
    9  = 742
    10 = 1047
    11 = 65
    12 = 329
    13 = 53
    14 = 6414
    15 = 12
    17 = 137
    18 = 152
    20 = 15355
    21 = 1523
    25 = 11049
    26 = 755
    27 = 4322
    28 = 239
    29 = 374
    30 = 18
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

# "全部类型" sheet has the same F-column values, but shifted down by one row
# starting at row 10 (because it has an extra data row compared to "展览")
$updates4 = @{
    3  = 528
    4  = 1526
    10 = 742
    12 = 1047
    13 = 65
    14 = 329
    15 = 53
    17 = 6414
    18 = 12
    20 = 137
    21 = 152
    23 = 15355
    24 = 1523
    28 = 11049
    29 = 755
    30 = 4322
    31 = 239
    32 = 374
    33 = 18
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
